$d = $word.ActiveDocument

# --- Paragraph 1: a new, truly empty paragraph after "3: return the result res." ---
$last = $d.Paragraphs.Last
$r = $last.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$p1 = $d.Paragraphs.Last
# Type a placeholder char so the paragraph mark's run-formatting is materialized,
# then delete just the character (not the paragraph mark) so no empty <w:r> remains.
$p1.Range.InsertAfter("X")
$p1b = $d.Paragraphs.Last
$sub1 = $d.Range($p1b.Range.Start, $p1b.Range.End - 1)
$sub1.Delete()

# --- Paragraph 2: "Do left - 1 so that I can keep the length bound but still be able to check 2 continuous value." ---
$p1c = $d.Paragraphs.Last
$r2 = $p1c.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$p2 = $d.Paragraphs.Last
$p2.Range.InsertAfter("Do left " + [char]0x2013 + " 1 so that I can keep the length bound but still be able to check 2 continuous value.")

# --- Paragraph 3: "Also it check the value 0" ---
$p2b = $d.Paragraphs.Last
$r3 = $p2b.Range
$r3.Collapse(0)
$r3.InsertParagraphAfter()

$p3 = $d.Paragraphs.Last
$p3.Range.InsertAfter("Also it check the value 0")

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
